$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.133.30"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.838.58"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'240.71"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").Value = "'0.6857"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'0.2990"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "'0.07414"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "'23.15"
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("D11").Value = "'0.07640"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "1.837.57"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "'5.052"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").Value = "'0.6805"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "'87.30"
$ws.Range("E15").Value = "  -6.57%  "
$ws.Range("D16").Value = "'6.150"
$ws.Range("E16").Value = "  -7.48%  "
$ws.Range("D17").Value = "29.137.35"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "'0.000008153"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").Value = "2.084.33"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "'229.22"
$ws.Range("E20").Value = "  -6.06%  "
$ws.Range("D21").Value = "'12.52"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'7.351"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'159.61"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.1440"
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("D27").Value = "'8.745"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "'18.07"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "'1.510"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'4.139"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "'0.05262"
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").Value = "'0.7546"
$ws.Range("E34").Value = "  -5.08%  "
$ws.Range("D35").Value = "'1.849"
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("D36").Value = "'1.132"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").Value = "'2.684"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "1.292.06"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").Value = "'0.01827"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("D40").Value = "'2.721"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").Value = "'0.9399"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'5.948"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("D43").Value = "'104.70"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "1.986.22"
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "'0.5196"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").Value = "'64.74"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("D49").Value = "'9.488"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").Value = "'1.765"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'0.07452"
$ws.Range("E51").Value = "  +17.98%  "
